$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.387.19"
$ws.Range("E2").Value = "  -2.99%  "

$ws.Range("D3").Value = "3.500.46"
$ws.Range("E3").Value = "  -4.71%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "603.43"
$ws.Range("E5").Value = "  -2.34%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "149.59"
$ws.Range("E6").Value = "  -6.11%  "

$ws.Range("D7").Value = "3.499.86"
$ws.Range("E7").Value = "  -4.62%  "

$ws.Range("E8").Value = "  -0.07%  "

$ws.Range("E9").Value = "  -3.29%  "

$ws.Range("E10").Value = "  -4.12%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.94"
$ws.Range("E11").Value = "  -3.59%  "

$ws.Range("E12").Value = "  -4.34%  "

$ws.Range("E13").Value = "  -4.34%  "

$ws.Range("D14").Value = "4.092.51"
$ws.Range("E14").Value = "  -4.61%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "31.47"
$ws.Range("E15").Value = "  -2.97%  "

$ws.Range("D16").Value = "3.505.12"
$ws.Range("E16").Value = "  -4.53%  "

$ws.Range("D17").Value = "67.267.07"
$ws.Range("E17").Value = "  -3.20%  "

$ws.Range("E18").Value = "  -0.68%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.36"
$ws.Range("E19").Value = "  -2.22%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.99"
$ws.Range("E20").Value = "  -5.67%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "446.34"
$ws.Range("E21").Value = "  -4.83%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.97"
$ws.Range("E22").Value = "  -12.77%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.618"
$ws.Range("E23").Value = "  -4.74%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "77.35"

$ws.Range("E25").Value = "  +5.50%  "

$ws.Range("E26").Value = "  +0.09%  "

$ws.Range("D27").Value = "3.640.03"
$ws.Range("E27").Value = "  -4.68%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.12"
$ws.Range("E28").Value = "  -8.76%  "

$ws.Range("E29").Value = "  -5.07%  "

$ws.Range("E30").Value = "  -4.92%  "

$ws.Range("E31").Value = "  -7.20%  "

$ws.Range("E32").Value = "  -0.01%  "

$ws.Range("E33").Value = "  +0.49%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.68"
$ws.Range("E34").Value = "  -3.50%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.10"
$ws.Range("E35").Value = "  -4.56%  "

$ws.Range("D36").Value = "3.488.34"
$ws.Range("E36").Value = "  -4.99%  "

$ws.Range("E38").Value = "  -3.89%  "

$ws.Range("E39").Value = "  +0.03%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  -0.03%  "

$ws.Range("E41").Value = "  -1.24%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "173.09"
$ws.Range("E42").Value = "  -3.11%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0874"
$ws.Range("E43").Value = "  -2.13%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.39"
$ws.Range("E44").Value = "  -7.01%  "

$ws.Range("E45").Value = "  -4.99%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "45.38"
$ws.Range("E46").Value = "  -3.01%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "27.68"
$ws.Range("E47").Value = "  -5.16%  "

$ws.Range("E48").Value = "  +5.48%  "

$ws.Range("E49").Value = "  -5.50%  "

$ws.Range("E50").Value = "  -4.23%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.997"
$ws.Range("E51").Value = "  -4.22%  "
